# Scheduled runner update: refresh cached Universalis market-board price
# snapshots (currentAveragePrice / currentAveragePriceNQ/HQ, LevePriceNQ/HQ
# and the derived LeveProfitNQ/HQ columns) across the per-class Leve profit
# sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1800
$ws.Range("I2").Value = 400
$ws.Range("K2").Value = 400
$ws.Range("M2").Value = -287
$ws.Range("H4").Value = 2793.625
$ws.Range("I4").Value = 2793.625
$ws.Range("K4").Value = 2793.625
$ws.Range("M4").Value = -2679.625
$ws.Range("H9").Value = 340.7143
$ws.Range("I9").Value = 512.5
$ws.Range("J9").Value = 111.666664
$ws.Range("K9").Value = 512.5
$ws.Range("L9").Value = 111.666664
$ws.Range("M9").Value = -343.5
$ws.Range("N9").Value = -449.666664
$ws.Range("H137").Value = 1773.7
$ws.Range("J137").Value = 2103.4
$ws.Range("L137").Value = 6310.200000000001
$ws.Range("N137").Value = -11410.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3596.3635
$ws.Range("I2").Value = 3656
$ws.Range("K2").Value = 3656
$ws.Range("M2").Value = -3543
$ws.Range("H5").Value = 136
$ws.Range("I5").Value = 95
$ws.Range("K5").Value = 95
$ws.Range("M5").Value = 17
$ws.Range("H21").Value = 216.66667
$ws.Range("I21").Value = 216.66667
$ws.Range("K21").Value = 216.66667
$ws.Range("M21").Value = 157.33333
$ws.Range("H61").Value = 2726.818
$ws.Range("I61").Value = 3077.4443
$ws.Range("J61").Value = 1149
$ws.Range("K61").Value = 3077.4443
$ws.Range("L61").Value = 1149
$ws.Range("M61").Value = -2865.4443
$ws.Range("N61").Value = -1573
$ws.Range("H74").Value = 1703.8235
$ws.Range("I74").Value = 1499.75
$ws.Range("K74").Value = 1499.75
$ws.Range("M74").Value = -625.75
$ws.Range("H77").Value = 1703.8235
$ws.Range("I77").Value = 1499.75
$ws.Range("K77").Value = 7498.75
$ws.Range("M77").Value = -3130.75
$ws.Range("H110").Value = 2238.8572
$ws.Range("I110").Value = 1767.7273
$ws.Range("K110").Value = 1767.7273
$ws.Range("M110").Value = 277.2727
$ws.Range("H116").Value = 3596.3635
$ws.Range("I116").Value = 3656
$ws.Range("K116").Value = 3656
$ws.Range("M116").Value = -1362
$ws.Range("H122").Value = 6591.731
$ws.Range("I122").Value = 6970.75
$ws.Range("J122").Value = 5328.3335
$ws.Range("K122").Value = 20912.25
$ws.Range("L122").Value = 15985.0005
$ws.Range("M122").Value = -18462.25
$ws.Range("N122").Value = -20885.0005
$ws.Range("H136").Value = 2726.818
$ws.Range("I136").Value = 3077.4443
$ws.Range("J136").Value = 1149
$ws.Range("K136").Value = 9232.332900000001
$ws.Range("L136").Value = 3447
$ws.Range("M136").Value = -6682.332900000001
$ws.Range("N136").Value = -8547

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3596.3635
$ws.Range("I3").Value = 3656
$ws.Range("K3").Value = 3656
$ws.Range("M3").Value = -3542
$ws.Range("H4").Value = 136
$ws.Range("I4").Value = 95
$ws.Range("K4").Value = 95
$ws.Range("M4").Value = 20
$ws.Range("H99").Value = 3671.353
$ws.Range("I99").Value = 4028.0667
$ws.Range("J99").Value = 996
$ws.Range("K99").Value = 4028.0667
$ws.Range("L99").Value = 996
$ws.Range("M99").Value = -2530.0667
$ws.Range("N99").Value = -3992
$ws.Range("H134").Value = 6949.8
$ws.Range("J134").Value = 9220
$ws.Range("L134").Value = 27660
$ws.Range("N134").Value = -32730

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6227.4287
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
$ws.Range("H34").Value = 6227.4287
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H58").Value = 2935.0908
$ws.Range("I58").Value = 3029.1
$ws.Range("J58").Value = 1995
$ws.Range("K58").Value = 3029.1
$ws.Range("L58").Value = 1995
$ws.Range("M58").Value = -2826.1
$ws.Range("N58").Value = -2401
$ws.Range("H132").Value = 2079.75
$ws.Range("I132").Value = 1322
$ws.Range("K132").Value = 3966
$ws.Range("M132").Value = -1436
$ws.Range("H134").Value = 2574.1
$ws.Range("I134").Value = 2530.875
$ws.Range("J134").Value = 2747
$ws.Range("K134").Value = 7592.625
$ws.Range("L134").Value = 8241
$ws.Range("M134").Value = -5057.625
$ws.Range("N134").Value = -13311
$ws.Range("H136").Value = 2935.0908
$ws.Range("I136").Value = 3029.1
$ws.Range("J136").Value = 1995
$ws.Range("K136").Value = 9087.299999999999
$ws.Range("L136").Value = 5985
$ws.Range("M136").Value = -6537.299999999999
$ws.Range("N136").Value = -11085

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2166262.5
$ws.Range("I4").Value = 2332871.5
$ws.Range("K4").Value = 6998614.5
$ws.Range("M4").Value = -6998502.5
$ws.Range("H8").Value = 1186.5
$ws.Range("I8").Value = 1186.5
$ws.Range("K8").Value = 3559.5
$ws.Range("M8").Value = -3420.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 4400400
$ws.Range("J21").Value = 2501000
$ws.Range("L21").Value = 2501000
$ws.Range("N21").Value = -2501346
$ws.Range("H30").Value = 4400400
$ws.Range("J30").Value = 2501000
$ws.Range("L30").Value = 2501000
$ws.Range("N30").Value = -2501210
$ws.Range("H36").Value = 10004000
$ws.Range("I36").Value = 20000000
$ws.Range("K36").Value = 20000000
$ws.Range("M36").Value = -19999515
$ws.Range("H102").Value = 3341
$ws.Range("I102").Value = 3267.2856
$ws.Range("K102").Value = 3267.2856
$ws.Range("M102").Value = -1645.2856
$ws.Range("H132").Value = 2985.4666
$ws.Range("I132").Value = 2985.4666
$ws.Range("K132").Value = 8956.399800000001
$ws.Range("M132").Value = -6426.399800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4383.353
$ws.Range("I7").Value = 4476.5
$ws.Range("J7").Value = 4159.8
$ws.Range("K7").Value = 4476.5
$ws.Range("L7").Value = 4159.8
$ws.Range("M7").Value = -4364.5
$ws.Range("N7").Value = -4383.8
$ws.Range("H16").Value = 5262.25
$ws.Range("I16").Value = 5262.25
$ws.Range("K16").Value = 5262.25
$ws.Range("M16").Value = -5092.25
$ws.Range("H126").Value = 4383.353
$ws.Range("I126").Value = 4476.5
$ws.Range("J126").Value = 4159.8
$ws.Range("K126").Value = 13429.5
$ws.Range("L126").Value = 12479.4
$ws.Range("M126").Value = -10959.5
$ws.Range("N126").Value = -17419.4
$ws.Range("H132").Value = 1752.5
$ws.Range("I132").Value = 1916.4166
$ws.Range("J132").Value = 769
$ws.Range("K132").Value = 5749.2498
$ws.Range("L132").Value = 2307
$ws.Range("M132").Value = -3219.2498
$ws.Range("N132").Value = -7367
$ws.Range("H136").Value = 3579.2727
$ws.Range("I136").Value = 4534
$ws.Range("K136").Value = 13602
$ws.Range("M136").Value = -11052

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1459.9584
$ws.Range("I136").Value = 1479.9565
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 4439.8695
$ws.Range("L136").Value = 3000
$ws.Range("M136").Value = -1889.8695
$ws.Range("N136").Value = -8100
